$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.869.73'
$ws.Range('E2').Value = '  +6.60%  '
$ws.Range('D3').Value = '2.307.79'
$ws.Range('E3').Value = '  +3.70%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.35'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.87'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +13.46%  '
$ws.Range('E7').Value = '  +2.69%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +8.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.01'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +13.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0805'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.44'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +7.20%  '
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('D14').Value = '2.659.84'
$ws.Range('E14').Value = '  +3.77%  '
$ws.Range('D15').Value = '2.301.24'
$ws.Range('E15').Value = '  +4.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.05'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.20%  '
$ws.Range('E17').Value = '  +5.80%  '
$ws.Range('D18').Value = '46.853.69'
$ws.Range('E18').Value = '  +6.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.88'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +26.49%  '
$ws.Range('D20').Value = '0.0₃0951'
$ws.Range('E20').Value = '  +5.41%  '
$ws.Range('E21').Value = '  +3.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '66.93'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.65'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +6.08%  '
$ws.Range('E24').Value = '  +5.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.96'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +6.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '43.75'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +15.35%  '
$ws.Range('E28').Value = '  +2.45%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.00'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +7.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.26'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +5.06%  '
$ws.Range('E31').Value = '  +14.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.83'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +8.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0808'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +7.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '147.78'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.22'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +14.67%  '
$ws.Range('E36').Value = '  +11.97%  '
$ws.Range('E37').Value = '  +3.48%  '
$ws.Range('E38').Value = '  +8.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.28'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +24.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.05'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +15.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.49'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +10.11%  '
$ws.Range('E42').Value = '  +2.42%  '
$ws.Range('E43').Value = '  +14.25%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = '1.858.06'
$ws.Range('E45').Value = '  +3.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.87'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +21.02%  '
$ws.Range('E47').Value = '  +10.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '74.91'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +12.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.93'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +11.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '97.41'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.84%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.42'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +8.13%  '
